{"js": "// Apply the commit: \"Added Abhishek's student number and changes in algorithm comparison\"\n//\n// 1. Author block: \"Abhishek Kakad\" -> \"Abhishek Kakad 1007545364\"\n// 2. Author block: \"Kevin Mano \" -> \"Kevin Mano 1007632992 \"\n// 3. Motivation section: \"that is able to learn\" -> \"that can learn\"\n// 4. Methodology section: \"...is formulated and solved by finding an H2 state\n//    feedback controller.\" -> \"...is formulated and an H2 state feedback\n//    controller is found by solving linear matrix inequalities.\"\n// 5. Comparisons bullet: drop the \"Unfortunately, ...\" sentence and replace it\n//    with \"Both MPC methods would not consider constraints.\"\n\n// 1) Abhishek Kakad -> append student number\nconst abhishek = context.document.body.search(\"Abhishek Kakad\", { matchCase: true });\nabhishek.load(\"text\");\nawait context.sync();\nif (abhishek.items.length > 0) {\n  abhishek.items[0].insertText(\" 1007545364\", Word.InsertLocation.end);\n}\n\n// 2) Kevin Mano -> insert student number right after the name (before the\n// trailing space that is already part of the paragraph).\nconst kevin = context.document.body.search(\"Kevin Mano\", { matchCase: true });\nkevin.load(\"text\");\nawait context.sync();\nif (kevin.items.length > 0) {\n  kevin.items[0].insertText(\" 1007632992\", Word.InsertLocation.after);\n}\n\n// 3) \"a Robust Controller that is able to learn\" -> \"... that can learn\"\nconst isAbleTo = context.document.body.search(\"is able to learn\", { matchCase: true });\nisAbleTo.load(\"text\");\nawait context.sync();\nif (isAbleTo.items.length > 0) {\n  isAbleTo.items[0].insertText(\"can learn\", Word.InsertLocation.replace);\n}\n\n// 4) Rework the ending of the methodology paragraph.\nconst methodOld = \"formulated and solved by finding an H2 state feedback controller.\";\nconst methodNew = \"formulated and an H2 state feedback controller is found by solving linear matrix inequalities.\";\nconst methodResults = context.document.body.search(methodOld, { matchCase: true });\nmethodResults.load(\"text\");\nawait context.sync();\nif (methodResults.items.length > 0) {\n  methodResults.items[0].insertText(methodNew, Word.InsertLocation.replace);\n}\n\n// 5) Replace the trailing sentence in the comparisons bullet.\nconst compareOld =\n  \" Unfortunately, we aren\\u2019t aware of a robust control implementation in the safe-control-gym environment. We are also open to feedback for algorithms to compare to generate fair comparisons.\";\nconst compareNew = \" Both MPC methods would not consider constraints.\";\nconst compareResults = context.document.body.search(compareOld, { matchCase: true });\ncompareResults.load(\"text\");\nawait context.sync();\nif (compareResults.items.length > 0) {\n  compareResults.items[0].insertText(compareNew, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Apply the commit: \"Added Abhishek's student number and changes in algorithm comparison\"\n#\n# 1. Author block: \"Abhishek Kakad\" -> \"Abhishek Kakad 1007545364\"\n# 2. Author block: \"Kevin Mano \" -> \"Kevin Mano 1007632992 \"\n# 3. Motivation section: \"that is able to learn\" -> \"that can learn\"\n# 4. Methodology section: \"...is formulated and solved by finding an H2 state\n#    feedback controller.\" -> \"...is formulated and an H2 state feedback\n#    controller is found by solving linear matrix inequalities.\"\n# 5. Comparisons bullet: drop the \"Unfortunately, ...\" sentence and replace it\n#    with \"Both MPC methods would not consider constraints.\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text {\n    param($SearchText, $NewText)\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Text = $SearchText\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Replacement.Text = $NewText\n    $range.Find.Forward = $true\n    $range.Find.Wrap = 0   # wdFindStop\n    $range.Find.MatchCase = $true\n    $range.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n}\n\n# 1) Abhishek Kakad -> append student number\nReplace-Text \"Abhishek Kakad\" \"Abhishek Kakad 1007545364\"\n\n# 2) Kevin Mano -> insert student number right after the name (keeping the\n# existing trailing space after it).\nReplace-Text \"Kevin Mano \" \"Kevin Mano 1007632992 \"\n\n# 3) \"a Robust Controller that is able to learn\" -> \"... that can learn\"\nReplace-Text \"is able to learn\" \"can learn\"\n\n# 4) Rework the ending of the methodology paragraph.\nReplace-Text \"formulated and solved by finding an H2 state feedback controller.\" \"formulated and an H2 state feedback controller is found by solving linear matrix inequalities.\"\n\n# 5) Replace the trailing sentence in the comparisons bullet.\n$rightQuote = [char]0x2019\n$oldCompareSentence = \" Unfortunately, we aren${rightQuote}t aware of a robust control implementation in the safe-control-gym environment. We are also open to feedback for algorithms to compare to generate fair comparisons.\"\nReplace-Text $oldCompareSentence \" Both MPC methods would not consider constraints.\"\n"}
